$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in Contest 23 (row 32) scores: RR vs DC ---
$ws.Range("E32").Value = 80
$ws.Range("H32").Value = 40
$ws.Range("K32").Value = 60
$ws.Range("N32").Value = 0
$ws.Range("Q32").Value = 20
$ws.Range("T32").Value = 100

# --- Insert a new blank contest row after row 41 (becomes row 42), pushing
#     the totals block down by one row. Copy row 41 (still blank at this
#     point, incl. formulas) into the new row so the RANK/VLOOKUP formulas
#     and borders carry over, mirroring the existing pattern used for every
#     contest row. ---
$ws.Range("A41:U41").Copy()
$ws.Range("A42:U42").Insert(-4121)
$excel.CutCopyMode = $false

# Re-copy formatting (and formulas) from row 41 into the new row 42,
# segment by segment so the blank separator columns (F, I, L, O, R) are
# left untouched.
$srcSegs = @("A41:C41", "D41:E41", "G41:H41", "J41:K41", "M41:N41", "P41:Q41", "S41:T41")
$dstSegs = @("A42:C42", "D42:E42", "G42:H42", "J42:K42", "M42:N42", "P42:Q42", "S42:T42")
for ($i = 0; $i -lt $srcSegs.Length; $i++) {
  $ws.Range($srcSegs[$i]).Copy()
  $ws.Range($dstSegs[$i]).PasteSpecial(-4122)
  $excel.CutCopyMode = $false
}

# --- Add Contest 32 (row 41, unaffected by the insert above since it shifted
#     everything starting at row 42 downward): MI vs KKR ---
$ws.Range("A41").Value = 32
$ws.Range("B41").Value = 1
$ws.Range("C41").Value = "MI vs KKR"

# --- Fix up the season Total row (now row 46) so its SUM ranges include
#     the newly inserted contest row 42. ---
$ws.Range("E46").Formula = "=SUM(D10:D42)"
$ws.Range("H46").Formula = "=SUM(G10:G42)"
$ws.Range("K46").Formula = "=SUM(J10:J42)"
$ws.Range("N46").Formula = "=SUM(M10:M42)"
$ws.Range("Q46").Formula = "=SUM(P10:P42)"
$ws.Range("T46").Formula = "=SUM(S10:S42)"

# --- The win/loss/tie conditional-formatting rules on the six "Total" cells
#     are anchored to explicit cells and don't follow a row insert on their
#     own, so re-point each rule group from row 45 to the new row 46. ---
$cfMoves = @(
  @("E45", "E46"),
  @("H45", "H46"),
  @("K45", "K46"),
  @("N45", "N46"),
  @("Q45", "Q46"),
  @("T45", "T46")
)
foreach ($mv in $cfMoves) {
  $fcs = $ws.Range($mv[0]).FormatConditions
  for ($i = 1; $i -le $fcs.Count; $i++) {
    [void]$fcs.Item($i).ModifyAppliesToRange($ws.Range($mv[1]))
  }
}

# Match the author's last-saved selection (row shifted by the insert above).
[void]$ws.Range("U46").Select()
